# docs/images: Update AddressBook to RestaurantBook in UndoRedo images
#
# The activity diagram on slide 1 has two shapes that mention the
# "address book" / "addressBookStateList" concept from the UndoRedo
# design. Rename them to "restaurant book" / "restaurantBookStateList"
# to match the updated RestaurantBook domain terminology, while leaving
# everything else (formatting, other runs, the rest of the diagram)
# untouched.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- "TextBox 47": "[command commits address book]" ------------------
# The whole second run ("command commits address book]") is replaced
# with "command commits restaurant book]"; the leading "[" run is left
# alone. This textbox auto-fits to its text (<a:spAutoFit/>), so pin
# its height back to the original after the edit (real PowerPoint
# re-lays autofit boxes out on every text change too).
$shpCommand = $s.Shapes.Item("TextBox 47")
$origHeight = $shpCommand.Height
$trCommand = $shpCommand.TextFrame.TextRange
$runCommand = $trCommand.Characters(2, $trCommand.Length - 1)
$runCommand.Text = "command commits restaurant book]"
$shpCommand.Height = $origHeight

# --- "Rectangle: Rounded Corners 50": "...save address book to ------
#      addressBookStateList "
# Only the variable-name run "addressBookStateList" changes to
# "restaurantBookStateList"; the sentence text before it (which still
# reads "...save address book to ...") and the trailing space run are
# left untouched.
$shpState = $s.Shapes.Item("Rounded Rectangle 50")
$trState = $shpState.TextFrame.TextRange
$fullState = $trState.Text
$startState = $fullState.IndexOf("addressBookStateList") + 1
$runState = $trState.Characters($startState, "addressBookStateList".Length)
$runState.Text = "restaurantBookStateList"
